$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "fullname"
$ws.Range("B1").Value = "phonenum"
$ws.Range("C1").Value = "cashapp"
$ws.Range("D1").Value = "numofspots"

$ws.Range("A2").Value = "Erik Bridges"
$ws.Range("B2").Value = " 555-555-555"
$ws.Range("C2").Value = ' $app'
$ws.Range("D2").Value = 5
